$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D26").Value = "인공지능 음성 생성 연구: 음성 분류 솔루션"

$ws.Range("D46").Value = "[CJ제일제당] 2022년 09월, 생물정보학(Bioinformatics 채용), R&D  직무 OMICS 연구원"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/485"

$ws.Range("D50").Value = "가장 쉬운 Self Organizing Maps (SOM) [Travelling Salesman Problem]"
$ws.Range("E50").Value = "http://incredible.egloos.com/7548275"

$wb.Save()
